$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.004031
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.034111333333333
$ws.Range("N2").Value = 3.102334
$ws.Range("O2").Value = 0.4722250593604582
$ws.Range("P2").Value = 0.4722250593604581
$ws.Range("Q2").Value = 0.004168502784666667
$ws.Range("R2").Value = 0.037516525062
$ws.Range("S2").Value = 0.4722250593604582
$ws.Range("T2").Value = 0.4722250593604581

# Row 3
$ws.Range("G3").Value = 0.004031
$ws.Range("M3").Value = 1.155758333333333
$ws.Range("N3").Value = 3.467275
$ws.Range("O3").Value = 0.5277749406395418
$ws.Range("P3").Value = 0.5277749406395419
$ws.Range("Q3").Value = 0.004658861841666666
$ws.Range("R3").Value = 0.041929756575
$ws.Range("S3").Value = 0.5277749406395418
$ws.Range("T3").Value = 0.5277749406395419
